$wb = $excel.ActiveWorkbook

# --- Sheet "Confirmed" (sheet1): add row 36 ---
$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsConfirmed.Range("A35:C35").Copy()
$wsConfirmed.Range("A36:C36").PasteSpecial(-4122)
$wsConfirmed.Range("A36").Value = 43932
$wsConfirmed.Range("B36").Formula = "=SUM(B35+C36)"
$wsConfirmed.Range("C36").Value = 58
$wsConfirmed.Range("B35:B36").Select()

# --- Sheet "Recoverd" (sheet2): add row 36 ---
$wsRecoverd = $wb.Worksheets.Item("Recoverd")
$wsRecoverd.Range("A35:C35").Copy()
$wsRecoverd.Range("A36:C36").PasteSpecial(-4122)
$wsRecoverd.Range("A36").Value = 43932
$wsRecoverd.Range("B36").Formula = "=SUM(B35+C36)"
$wsRecoverd.Range("C36").Value = 3
$wsRecoverd.Range("C39").Select()

# --- Sheet "Death" (sheet3): add row 36 ---
$wsDeath = $wb.Worksheets.Item("Death")
$wsDeath.Range("A35:C35").Copy()
$wsDeath.Range("A36:C36").PasteSpecial(-4122)
$wsDeath.Range("A36").Value = 43932
$wsDeath.Range("B36").Formula = "=SUM(B35+C36)"
$wsDeath.Range("C36").Value = 3
$wsDeath.Range("C38").Select()

# Death becomes the active/selected tab (matches activeTab="2" and tabSelected on sheet3)
$wsDeath.Activate()
